# Refactored check_plea_and_findings from charges grid class.
# Appends the new charge rows (64-69) that resulted from that refactor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row, ColA..ColI. Values that look numeric ("0", "4510.11", "4511.20") are
# prefixed with a leading apostrophe so Excel stores them as text, matching
# the other charge rows in this column (case numbers / statute codes / the
# "0" flag column are all text in this sheet).
$rows = @(
    @{ Row=64; A="21CRB01291"; B="Bunner"; C="PERMISSION REQ'D TO USE LICENSED DOCK"; D="1501:46-12-04";  E="MM";                 F="Guilty";     G="Guilty"; H=0; I="'0" },
    @{ Row=65; A="21CRB01291"; B="Bunner"; C="PERMISSION REQ'D TO USE LICENSED DOCK"; D="1501:46-12-04";  E="MM";                 F="No Contest"; G="Guilty"; H=0; I="'0" },
    @{ Row=66; A="21CRB01291"; B="Bunner"; C="No Operator License - Expired";         D="4510.12(C)(2)"; E="Minor Misdemeanor";  F="No Contest"; G="Guilty"; H=0; I="'0" },
    @{ Row=67; A="21TRD09437"; B="Bunner"; C="DUS";                                   D="'4510.11";       E="M1";                 F="Guilty";     G="Guilty"; H=0; I="'0" },
    @{ Row=68; A="21TRD09437"; B="Bunner"; C="1ST SPEED 1 YR SCHOOL >35MPHM4";        D="4511.21B1A";     E="M4";                 F="Guilty";     G="Guilty"; H=0; I="'0" },
    @{ Row=69; A="21TRD09437"; B="Bunner"; C="RECKLESS OPERATION 1ST IN 1 YR";        D="'4511.20";       E="MM";                 F="Dismissed";  G="";        H=0; I="'0" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
    $ws.Cells.Item($n, 9).Value = $r.I
}
